$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Create WO")

# Update cell B2's text value
$ws.Range("B2").Value = "Pro-SYDATA1 (Lot track)"

# Move the active selection to B2 (single-cell selection instead of whole row)
$ws.Range("B2").Select()
